# Auto commit at 2025-08-30  8:06:05.21
# Update the Metrics sheet figures (monthly refresh of the underlying data)
# and move the "active" sheet/selection from Metrics to the today sheet.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# --- Refresh the Metrics values (B2:B13) -----------------------------------
# The "today" sheet pulls these in live via =Metrics!Bn formulas, so updating
# them here also refreshes today!B11:B22, E11:E22 and F11:F22 automatically.
$metrics.Range("B2").Value  = 489786.06
$metrics.Range("B3").Value  = 420294.83
$metrics.Range("B4").Value  = 155041.43
$metrics.Range("B5").Value  = 19203
$metrics.Range("B6").Value  = 3885414.6299999994
$metrics.Range("B7").Value  = 3299009.4899999998
$metrics.Range("B8").Value  = 1118683.99
$metrics.Range("B9").Value  = 149891
$metrics.Range("B10").Value = 32350738.43099983
$metrics.Range("B11").Value = 19328879.560000002
$metrics.Range("B12").Value = 11400392.880000001
$metrics.Range("B13").Value = 1247518

# --- Move the active sheet/selection ---------------------------------------
# Previously Metrics was the active tab with E19 selected; now the today
# sheet is active (workbook activeTab=3) with H8 selected, and Metrics keeps
# a plain (non-active) selection of G41.
$metrics.Activate()
$metrics.Range("G41").Select()

$today.Activate()
$today.Range("H8").Select()
